# Add a second patient record (row 3) to the "исходный формат" sheet,
# mirroring the layout/formatting already used by row 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Carry the number formats (date/time styles) from row 2 down to row 3
# before writing values, so the new row re-uses the existing date/time
# cell styles instead of Excel minting brand-new ones.
$ws.Range("B2:N2").Copy()
$ws.Range("B3:N3").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("A3").Value = "Петров Пётр Петрович"
$ws.Range("B3").Value = 114303
$ws.Range("C3").Value = 34
$ws.Range("D3").Value = "Пушкина 12"
$ws.Range("E3").Value = "В23-32"
$ws.Range("F3").Value = "Склероз"
$ws.Range("G3").Value = "Обыкновенное"
$ws.Range("H3").Value = "Пушкина 13"
$ws.Range("I3").Value = "Лучшая лаборатория мира"
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 114304
$ws.Range("L3").Value = 0.4236111111111111
$ws.Range("M3").Value = 114305
$ws.Range("N3").Value = 0.95833333333333337

# New "Название лаборатории" column (I) needs its width best-fit to the
# newly added content, like the other text columns on this sheet.
$ws.Columns.Item(9).AutoFit()

# Match the cursor position left behind in the saved file.
$ws.Range("G10").Select() | Out-Null
